# Update "Energy Consumption15" worksheet data with new computed values
# and extend the table with one additional row (A49:C49), per the
# "using priority queue, improved performance" re-computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2; A = 0; B = 7.978418581042527; C = 0.7408803578404031 },
    @{ Row = 3; A = 1; B = 14.47477963335235; C = 1.4896011090474 },
    @{ Row = 4; A = 2; B = 14.59153928952939; C = 2.325090712013005 },
    @{ Row = 5; A = 3; B = 17.08897116917147; C = 3.206358666379736 },
    @{ Row = 6; A = 4; B = 18.42736315900551; C = 4.119545177541847 },
    @{ Row = 7; A = 5; B = 18.52899791702236; C = 4.916245899095418 },
    @{ Row = 8; A = 6; B = 24.65937219061949; C = 5.727784297528979 },
    @{ Row = 9; A = 7; B = 25.87882233437227; C = 6.572569400737355 },
    @{ Row = 10; A = 8; B = 29.02464985940152; C = 7.45878303424291 },
    @{ Row = 11; A = 9; B = 29.11796794180655; C = 8.353895251220521 },
    @{ Row = 12; A = 10; B = 29.87496845137489; C = 9.575642558799883 },
    @{ Row = 13; A = 11; B = 29.99778393585725; C = 10.54000525353741 },
    @{ Row = 14; A = 12; B = 37.08746538175033; C = 11.70490030839822 },
    @{ Row = 15; A = 13; B = 37.89073514459027; C = 12.62902853606904 },
    @{ Row = 16; A = 14; B = 37.95772166562909; C = 13.52938127026607 },
    @{ Row = 17; A = 15; B = 38.87991126865267; C = 14.2734350313667 },
    @{ Row = 18; A = 16; B = 39.26863240922653; C = 15.23408444296631 },
    @{ Row = 19; A = 17; B = 39.38645467148035; C = 16.26981379724187 },
    @{ Row = 20; A = 18; B = 41.90646476217339; C = 17.08572752911935 },
    @{ Row = 21; A = 19; B = 42.15913665555601; C = 17.95484496064131 },
    @{ Row = 22; A = 20; B = 42.19407510612744; C = 18.82954127385875 },
    @{ Row = 23; A = 21; B = 42.92703715619879; C = 19.84651201071057 },
    @{ Row = 24; A = 22; B = 43.94464527089828; C = 20.71802428192277 },
    @{ Row = 25; A = 23; B = 44.2724540647949; C = 21.592662500572 },
    @{ Row = 26; A = 24; B = 45.14361775815542; C = 22.32240559561387 },
    @{ Row = 27; A = 25; B = 55.49259450074115; C = 23.26361592395075 },
    @{ Row = 28; A = 26; B = 60.09635938348296; C = 24.14476736346783 },
    @{ Row = 29; A = 27; B = 61.47530979211123; C = 24.95872323212815 },
    @{ Row = 30; A = 28; B = 61.53804947706858; C = 25.84761512960964 },
    @{ Row = 31; A = 29; B = 61.66730947973826; C = 26.83911509348242 },
    @{ Row = 32; A = 30; B = 63.33661129260693; C = 27.70889178186824 },
    @{ Row = 33; A = 31; B = 66.18426747375798; C = 28.56221335582747 },
    @{ Row = 34; A = 32; B = 66.25355162560798; C = 29.5123779748203 },
    @{ Row = 35; A = 33; B = 69.41584862319266; C = 30.31934161777221 },
    @{ Row = 36; A = 34; B = 69.89414603213366; C = 31.33460733835416 },
    @{ Row = 37; A = 35; B = 69.96897584597671; C = 32.27800753011753 },
    @{ Row = 38; A = 36; B = 75.12501149406354; C = 33.16626417228397 },
    @{ Row = 39; A = 37; B = 75.22537377423372; C = 34.10957171608733 },
    @{ Row = 40; A = 38; B = 77.03067058810265; C = 34.98768375666029 },
    @{ Row = 41; A = 39; B = 77.08711983727127; C = 35.99891035389001 },
    @{ Row = 42; A = 40; B = 77.4402315839798; C = 36.9384740180268 },
    @{ Row = 43; A = 41; B = 79.77165615589128; C = 37.88935953937701 },
    @{ Row = 44; A = 42; B = 79.84733342660105; C = 38.75037202071826 },
    @{ Row = 45; A = 43; B = 82.41453877158258; C = 40.0533447745875 },
    @{ Row = 46; A = 44; B = 82.57651182041016; C = 40.91375327094882 },
    @{ Row = 47; A = 45; B = 92.34659287029028; C = 41.74620013022608 },
    @{ Row = 48; A = 46; B = 92.60165305046408; C = 42.66165393958249 },
    @{ Row = 49; A = 47; B = 93.58036577957165; C = 43.55683382436312 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
}

# Ensure the new row's A cell carries the same formatting as the rest of
# column A (border + centered/top alignment), matching the existing rows.
$ws.Range("A48").Copy()
$ws.Range("A49").PasteSpecial(-4122)
